# TC19_Canine_Filter_Diagnosis-RespCarciStg4.xlsx
# Insert a new "StatQuery" column between the existing "query" (A) and
# "dbExcel" (B) columns, shifting dbExcel -> C and WebExcel -> D.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at B; everything from the old B onward shifts right.
$ws.Columns.Item(2).Insert()

# New header + row-2 value for the inserted "StatQuery" column.
$ws.Range("B1").Value = "StatQuery"
$ws.Range("B2").Value = 'MATCH (s:study) WITH COLLECT(DISTINCT(s.clinical_study_designation)) AS all_studies MATCH (d:demographic) WITH COLLECT(DISTINCT(d.breed)) AS all_breeds, COLLECT(DISTINCT(d.sex)) AS all_sexes, all_studies MATCH (d:diagnosis) WITH COLLECT(DISTINCT(d.disease_term)) AS all_diseases, all_breeds, all_sexes, all_studies MATCH (p:program)<-[*]-(s:study)<-[*]-(c:case)<--(demo:demographic), (c)<--(diag:diagnosis) WHERE diag.disease_term IN[''Stage 4'']  OPTIONAL MATCH (f:file)-[*]->(c), (samp:sample)-[*]->(c) WITH DISTINCT c AS c, p, s, demo, diag, f, samp RETURN count(DISTINCT(f)) as number_of_files , count(DISTINCT(samp)) as number_of_sample , count(DISTINCT(c.case_id)) as number_of_cases , count(DISTINCT(s.clinical_study_designation)) as number_of_study'

# Match column B's width to column A's (~75.8) and wrap the long query text.
$ws.Columns.Item(2).ColumnWidth = 75
$ws.Range("B2").WrapText = $true

# Selection moves to A2 (matching the post-edit view state).
$ws.Range("A2").Select()
